$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.797.31'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '1.917.21'
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("E5").Value = '  -2.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4915'
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2976'
$ws.Range("E8").Value = '  +0.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06770'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").Value = '1.884.95'
$ws.Range("E10").Value = '  -0.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.10'
$ws.Range("E11").Value = '  -0.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07317'
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.187'
$ws.Range("E13").Value = '  +2.55%  '

$ws.Range("E14").Value = '  -1.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6729'
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").Value = '30.784.89'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007985'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.57'
$ws.Range("E18").Value = '  +2.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").Value = '2.138.06'
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.280'
$ws.Range("E22").Value = '  +9.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '203.18'
$ws.Range("E23").Value = '  +10.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.272'
$ws.Range("E24").Value = '  +3.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.657'
$ws.Range("E25").Value = '  +3.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.07'
$ws.Range("E26").Value = '  +3.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.94'
$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.971'
$ws.Range("E28").Value = '  +3.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.435'
$ws.Range("E29").Value = '  +2.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.353'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09213'
$ws.Range("E31").Value = '  +2.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05717'
$ws.Range("E32").Value = '  +10.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.080'
$ws.Range("E33").Value = '  +2.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7502'
$ws.Range("E34").Value = '  +1.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.122'
$ws.Range("E35").Value = '  +1.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.705'
$ws.Range("E36").Value = '  -1.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01863'
$ws.Range("E37").Value = '  +1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.724'
$ws.Range("E38").Value = '  +2.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9270'
$ws.Range("E39").Value = '  -1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.081'
$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4506'
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.53'
$ws.Range("E42").Value = '  +25.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.59'
$ws.Range("E43").Value = '  +2.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.923'
$ws.Range("E44").Value = '  +2.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.009'
$ws.Range("E45").Value = '  +0.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1393'
$ws.Range("E46").Value = '  +4.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.730'
$ws.Range("E47").Value = '  +1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.09'
$ws.Range("E48").Value = '  +8.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06047'
$ws.Range("E49").Value = '  +3.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.059'
$ws.Range("E50").Value = '  +4.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4066'
$ws.Range("E51").Value = '  +3.67%  '
